$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = [double]"499.9999999999999"
$ws.Range("C3").Value2 = [double]"111.3217836438596"
$ws.Range("C4").Value2 = [double]"52469.16847303212"
$ws.Range("C5").Value2 = [double]"49929.83435027857"
$ws.Range("C6").Value2 = [double]"2673.1684730321313"
$ws.Range("C7").Value2 = [double]"2539.510049380525"
$ws.Range("C10").Value2 = [double]"130.0"
$ws.Range("C12").Value2 = [double]"0.5343075876139025"
$ws.Range("C13").Value2 = [double]"97.6809609675536"
$ws.Range("C14").Value2 = [double]"141.70577064630237"
$ws.Range("C15").Value2 = [double]"95.44884708572303"
$ws.Range("C16").Value2 = [double]"22.09108979931446"
$ws.Range("C17").Value2 = [double]"99.99999999999994"
$ws.Range("C18").Value2 = [double]"38.39482262187994"
$ws.Range("C20").Value2 = [double]"3.586709101989584"
$ws.Range("C21").Value2 = [double]"0.5574276500924213"
$ws.Range("C23").Value2 = [double]"0.404573804951899"
$ws.Range("C24").Value2 = [double]"16.58268925796842"
$ws.Range("C25").Value2 = [double]"17.534639693169538"
$ws.Range("C26").Value2 = [double]"18.999999999999993"
$ws.Range("C27").Value2 = [double]"4.491605370931133"
$ws.Range("C28").Value2 = [double]"12.990364158819158"
$ws.Range("C31").Value2 = [double]"0.9666666666666615"
$ws.Range("C32").Value2 = [double]"0.3512446913528251"
$ws.Range("C34").Value2 = [double]"37.5384152867733"
$ws.Range("C35").Value2 = [double]"797.0916524655247"
$ws.Range("C36").Value2 = [double]"616.6599238333979"
$ws.Range("C37").Value2 = [double]"1.5514929728291236"
$ws.Range("C38").Value2 = [double]"270.88704679019213"
$ws.Range("C39").Value2 = [double]"671.6959858942153"
$ws.Range("C40").Value2 = [double]"21.57065020054756"
$ws.Range("C41").Value2 = [double]"122.54649450899547"
$ws.Range("C42").Value2 = [double]"-0.031612571950745405"
$ws.Range("C43").Value2 = [double]"-2.2737367544323206E-13"
$ws.Range("C45").Value2 = [double]"52431.80598437237"
$ws.Range("C46").Value2 = [double]"51634.71433190684"
$ws.Range("C47").Value2 = [double]"51018.05440807344"
$ws.Range("C48").Value2 = [double]"51016.5029151006"
$ws.Range("C49").Value2 = [double]"50745.61586831041"
$ws.Range("C50").Value2 = [double]"50073.91988241619"
$ws.Range("C51").Value2 = [double]"50052.34923221564"
$ws.Range("C52").Value2 = [double]"49929.80273770663"
$ws.Range("C53").Value2 = [double]"49929.83435027858"
$ws.Range("C54").Value2 = [double]"49929.83435027857"
$ws.Range("C58").Value2 = [double]"147.2287526559152"
$ws.Range("C60").Value2 = [double]"0.22260203694624983"
$ws.Range("C62").Value2 = [double]"2.126189561697079"
$ws.Range("C64").Value2 = [double]"0.0926496353807621"
$ws.Range("C66").Value2 = [double]"22.948709435920396"
$ws.Range("C67").Value2 = [double]"46611.949108155946"
$ws.Range("C68").Value2 = [double]"37395.09870081817"
$ws.Range("C70").Value2 = [double]"7036.575356478633"
$ws.Range("C74").Value2 = [double]"424.73295727057376"
$ws.Range("C76").Value2 = [double]"0.7205798230090309"
$ws.Range("C77").Value2 = [double]"0.4154682889938561"
$ws.Range("C78").Value2 = [double]"1.091867864657074"
$ws.Range("C79").Value2 = [double]"0.029624798604174715"
$ws.Range("C80").Value2 = [double]"0.0737133924294191"
$ws.Range("C81").Value2 = [double]"14.024341381862035"
$ws.Range("C82").Value2 = [double]"14.812340453636596"
$ws.Range("C83").Value2 = [double]"24089.471436028558"
$ws.Range("C84").Value2 = [double]"11058.292826177614"
$ws.Range("C85").Value2 = [double]"8242.460890642398"
$ws.Range("C86").Value2 = [double]"7850.13349417967"
$ws.Range("C89").Value2 = [double]"484.28579957027137"
$ws.Range("C90").Value2 = [double]"485.66476888215806"
$ws.Range("C91").Value2 = [double]"0.8216140748357936"
$ws.Range("C92").Value2 = [double]"0.8239535623789305"
$ws.Range("C93").Value2 = [double]"0.31471200573089897"
$ws.Range("C94").Value2 = [double]"0.307316933319692"
$ws.Range("C95").Value2 = [double]"0.024749736221345076"
$ws.Range("C96").Value2 = [double]"0.024576978216942744"
$ws.Range("C97").Value2 = [double]"12.71577211637027"
$ws.Range("C98").Value2 = [double]"12.504260312516188"
$ws.Range("C99").Value2 = [double]"8952.272677352392"
$ws.Range("C100").Value2 = [double]"8940.482191946496"
$ws.Range("C101").Value2 = [double]"8952.272677352392"
$ws.Range("C102").Value2 = [double]"8940.482191946496"
$ws.Range("C105").Value2 = [double]"359.38939993419694"
$ws.Range("C106").Value2 = [double]"250.550253192183"
$ws.Range("C107").Value2 = [double]"0.6097213455251858"
$ws.Range("C108").Value2 = [double]"0.3807413742403091"
$ws.Range("C109").Value2 = [double]"0.5641559818275375"
$ws.Range("C110").Value2 = [double]"0.5638119644814734"
$ws.Range("C111").Value2 = [double]"0.032532257073064265"
$ws.Range("C112").Value2 = [double]"0.032518123081162316"
$ws.Range("C113").Value2 = [double]"17.341433782491585"
$ws.Range("C114").Value2 = [double]"17.338391981426767"
$ws.Range("C115").Value2 = [double]"-311.0648125400993"
$ws.Range("C116").Value2 = [double]"808.723437359591"
$ws.Range("C117").Value2 = [double]"7132.1423557438675"
$ws.Range("C118").Value2 = [double]"7129.035804240933"
$ws.Range("C125").Value2 = [double]"0.40425347679986806"
$ws.Range("C126").Value2 = [double]"0.6112437681927804"
$ws.Range("C127").Value2 = [double]"0.02920493943370731"
$ws.Range("C128").Value2 = [double]"0.03821842225421682"
$ws.Range("C129").Value2 = [double]"13.841955663612607"
$ws.Range("C130").Value2 = [double]"15.993432803870885"
$ws.Range("C133").Value2 = [double]"7956.682571278696"
$ws.Range("C134").Value2 = [double]"6982.0789150678"
$ws.Range("C137").Value2 = [double]"461.6611808872903"
$ws.Range("C138").Value2 = [double]"461.6611808872903"
$ws.Range("C139").Value2 = [double]"0.7369292919263221"
$ws.Range("C140").Value2 = [double]"0.7369292919263221"
$ws.Range("C141").Value2 = [double]"0.2026162542837819"
$ws.Range("C142").Value2 = [double]"0.19821489017607719"
$ws.Range("C143").Value2 = [double]"0.02256092885805613"
$ws.Range("C144").Value2 = [double]"0.022497079317791584"
$ws.Range("C145").Value2 = [double]"8.98084717870253"
$ws.Range("C146").Value2 = [double]"8.810694373972401"
$ws.Range("C147").Value2 = [double]"12457.057833945677"
$ws.Range("C148").Value2 = [double]"12445.464391336833"
$ws.Range("C149").Value2 = [double]"12457.057833945677"
$ws.Range("C150").Value2 = [double]"12445.464391336833"
$ws.Range("C157").Value2 = [double]"0.5533954471594784"
$ws.Range("C158").Value2 = [double]"0.5526559730475576"
$ws.Range("C159").Value2 = [double]"0.03209772622566322"
$ws.Range("C160").Value2 = [double]"0.03207062278946487"
$ws.Range("C161").Value2 = [double]"17.240954803739964"
$ws.Range("C162").Value2 = [double]"17.232467753295516"
$ws.Range("C163").Value2 = [double]"1130.1946133552606"
$ws.Range("C164").Value2 = [double]"445.6387321062796"
$ws.Range("C165").Value2 = [double]"7034.972385892438"
$ws.Range("C166").Value2 = [double]"7028.294773731165"
$ws.Range("C173").Value2 = [double]"0.7189555552356247"
$ws.Range("C174").Value2 = [double]"0.7176346188878386"
$ws.Range("C175").Value2 = [double]"0.03944756059120251"
$ws.Range("C176").Value2 = [double]"0.03938157421584426"
$ws.Range("C177").Value2 = [double]"18.225602406349157"
$ws.Range("C178").Value2 = [double]"18.222598592798633"
$ws.Range("C179").Value2 = [double]"6054.479788035808"
$ws.Range("C180").Value2 = [double]"6044.3520850827035"
$ws.Range("C181").Value2 = [double]"6054.479788035808"
$ws.Range("C182").Value2 = [double]"6044.3520850827035"
$ws.Range("C189").Value2 = [double]"0.5513919548920491"
$ws.Range("C190").Value2 = [double]"0.5513542365450854"
$ws.Range("C191").Value2 = [double]"0.032024293601087474"
$ws.Range("C192").Value2 = [double]"0.032022911136449485"
$ws.Range("C193").Value2 = [double]"17.217927169932175"
$ws.Range("C194").Value2 = [double]"17.21749263193991"
$ws.Range("C195").Value2 = [double]"-54.63182206511402"
$ws.Range("C196").Value2 = [double]"13.13648857717249"
$ws.Range("C197").Value2 = [double]"7016.880414700567"
$ws.Range("C198").Value2 = [double]"7016.539809818896"
$ws.Range("C201").Value2 = [double]"117.31576827160768"
$ws.Range("C203").Value2 = [double]"0.17738434640210213"
$ws.Range("C213").Value2 = [double]"4559.668792783809"
